# Revert "Revert "Additional feature:""
#
# 1. Rename the sheet from "SAP" back to "grs".
# 2. Re-apply rich (mixed-run) formatting to the "ACTIONS TAKEN" column
#    header (I8): the parenthetical instructions become italic, and the
#    closing parenthesis is explicitly bold (matching the already-bold
#    cell font), producing a multi-run shared string instead of a single
#    plain-text run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Sheet name -----------------------------------------------------
$ws.Name = "grs"

# --- 2. Rich text on the "ACTIONS TAKEN..." header cell -----------------
$cell = $ws.Cells.Item(8, 9)   # I8

$fullText = $cell.Text
$italicPhrase = "If referred, please specify name to whom grievance was referred"
$closeParen = ")"

$italicStart = $fullText.IndexOf($italicPhrase) + 1          # 1-based index
$italicLength = $italicPhrase.Length
$parenStart = $fullText.IndexOf($closeParen, $italicStart + $italicLength - 1) + 1
$parenLength = $closeParen.Length

# Middle run: "If referred, please specify name to whom grievance was
# referred" -> italic, not bold, Arial 11.
$middleRun = $cell.Characters($italicStart, $italicLength)
$middleRun.Font.Name = "Arial"
$middleRun.Font.Size = 11
$middleRun.Font.Bold = $false
$middleRun.Font.Italic = $true

# Trailing run: ")" -> keep bold (same as the rest of the header), Arial 11.
$trailingRun = $cell.Characters($parenStart, $parenLength)
$trailingRun.Font.Name = "Arial"
$trailingRun.Font.Size = 11
$trailingRun.Font.Bold = $true
